$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet regardless of its current name
$ws.Name = "Superdetails"

# Reset the stored selection back to the default top-left cell
$ws.Range("A1").Select()
